$d = $word.ActiveDocument

# Locate the paragraph that holds the "m:'A sample table'.bothMerge()" field
# and rewrite its field-code runs (fldChar/instrText) into plain literal-text
# runs, e.g. {m:'A sample table'.bothMerge()} - i.e. "split" the field into
# plain template tokens while keeping the _GoBack bookmark in place.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $target = $p
    }
}

$r = $target.Range

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F"><w:r><w:t>{</w:t></w:r><w:r w:rsidR="00DE6D5A"><w:t>m</w:t></w:r><w:r w:rsidR="00CB78EF"><w:t>:'A</w:t></w:r><w:r w:rsidR="001C5C89"><w:t xml:space="preserve"> sample</w:t></w:r><w:r w:rsidR="00CB78EF"><w:t xml:space="preserve"> table'.</w:t></w:r><w:r w:rsidR="00AC480E"><w:t>both</w:t></w:r><w:r w:rsidR="00474E6A"><w:t>Merge</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="00CB78EF"><w:t>()</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r.InsertXML($xml)
